# Auto-generated edit script: applies Universalis market-price refresh
# values (currentAveragePrice* / LevePrice* / LeveProfit* columns) across
# all eight Leve-profit tables (one per crafting job).
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H28").Value = 3560.5264
$ws_ALC.Range("I28").Value = 2182
$ws_ALC.Range("K28").Value = 2182
$ws_ALC.Range("M28").Value = -1697
$ws_ALC.Range("H76").Value = 9829.576999999999
$ws_ALC.Range("J76").Value = 8987.5
$ws_ALC.Range("L76").Value = 8987.5
$ws_ALC.Range("N76").Value = -9617.5
$ws_ALC.Range("H79").Value = 9829.576999999999
$ws_ALC.Range("J79").Value = 8987.5
$ws_ALC.Range("L79").Value = 8987.5
$ws_ALC.Range("N79").Value = -11171.5
$ws_ALC.Range("H98").Value = 3305.8604
$ws_ALC.Range("I98").Value = 3344.7354
$ws_ALC.Range("K98").Value = 3344.7354
$ws_ALC.Range("M98").Value = -1846.7354
$ws_ALC.Range("H122").Value = 3305.8604
$ws_ALC.Range("I122").Value = 3344.7354
$ws_ALC.Range("K122").Value = 10034.2062
$ws_ALC.Range("M122").Value = -7584.206200000001
$ws_ALC.Range("H132").Value = 1372.1818
$ws_ALC.Range("I132").Value = 1444.359
$ws_ALC.Range("K132").Value = 4333.076999999999
$ws_ALC.Range("M132").Value = -1803.076999999999
$ws_ALC.Range("H137").Value = 2714.3635
$ws_ALC.Range("I137").Value = 2416.625
$ws_ALC.Range("J137").Value = 2884.5
$ws_ALC.Range("K137").Value = 7249.875
$ws_ALC.Range("L137").Value = 8653.5
$ws_ALC.Range("M137").Value = -4699.875
$ws_ALC.Range("N137").Value = -13753.5
$ws_ALC.Range("H138").Value = 2131860.8
$ws_ALC.Range("I138").Value = 2760.8462
$ws_ALC.Range("J138").Value = 2945928.5
$ws_ALC.Range("K138").Value = 8282.5386
$ws_ALC.Range("L138").Value = 8837785.5
$ws_ALC.Range("M138").Value = -3142.5386
$ws_ALC.Range("N138").Value = -8848065.5

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H74").Value = 46240.293
$ws_ARM.Range("I74").Value = 84861.914
$ws_ARM.Range("J74").Value = 7618.6665
$ws_ARM.Range("K74").Value = 84861.914
$ws_ARM.Range("L74").Value = 7618.6665
$ws_ARM.Range("M74").Value = -83987.914
$ws_ARM.Range("N74").Value = -9366.666499999999
$ws_ARM.Range("H77").Value = 46240.293
$ws_ARM.Range("I77").Value = 84861.914
$ws_ARM.Range("J77").Value = 7618.6665
$ws_ARM.Range("K77").Value = 424309.57
$ws_ARM.Range("L77").Value = 38093.3325
$ws_ARM.Range("M77").Value = -419941.57
$ws_ARM.Range("N77").Value = -46829.3325
$ws_ARM.Range("H102").Value = 1272.4
$ws_ARM.Range("I102").Value = 1283.75
$ws_ARM.Range("K102").Value = 1283.75
$ws_ARM.Range("M102").Value = 338.25
$ws_ARM.Range("H132").Value = 750314.4
$ws_ARM.Range("I132").Value = 1139064.9
$ws_ARM.Range("J132").Value = 6617.7827
$ws_ARM.Range("K132").Value = 3417194.7
$ws_ARM.Range("L132").Value = 19853.3481
$ws_ARM.Range("M132").Value = -3414664.7
$ws_ARM.Range("N132").Value = -24913.3481

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H99").Value = 5351446
$ws_BSM.Range("I99").Value = 3191.7693
$ws_BSM.Range("K99").Value = 3191.7693
$ws_BSM.Range("M99").Value = -1693.7693
$ws_BSM.Range("H124").Value = 51507
$ws_BSM.Range("J124").Value = 51507
$ws_BSM.Range("L124").Value = 51507
$ws_BSM.Range("N124").Value = -61327

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H16").Value = 3580.4062
$ws_CRP.Range("I16").Value = 3048.96
$ws_CRP.Range("J16").Value = 5478.4287
$ws_CRP.Range("K16").Value = 3048.96
$ws_CRP.Range("L16").Value = 5478.4287
$ws_CRP.Range("M16").Value = -2761.96
$ws_CRP.Range("N16").Value = -6052.4287
$ws_CRP.Range("H31").Value = 5609.787
$ws_CRP.Range("I31").Value = 3000.6562
$ws_CRP.Range("J31").Value = 8488.826999999999
$ws_CRP.Range("K31").Value = 3000.6562
$ws_CRP.Range("L31").Value = 8488.826999999999
$ws_CRP.Range("M31").Value = -2705.6562
$ws_CRP.Range("N31").Value = -9078.826999999999
$ws_CRP.Range("H34").Value = 5609.787
$ws_CRP.Range("I34").Value = 3000.6562
$ws_CRP.Range("J34").Value = 8488.826999999999
$ws_CRP.Range("K34").Value = 3000.6562
$ws_CRP.Range("L34").Value = 8488.826999999999
$ws_CRP.Range("M34").Value = -2798.6562
$ws_CRP.Range("N34").Value = -8892.826999999999
$ws_CRP.Range("H52").Value = 59500
$ws_CRP.Range("J52").Value = 59500
$ws_CRP.Range("L52").Value = 59500
$ws_CRP.Range("N52").Value = -60088
$ws_CRP.Range("H107").Value = 1222.1818
$ws_CRP.Range("I107").Value = 639.86664
$ws_CRP.Range("J107").Value = 2470
$ws_CRP.Range("K107").Value = 639.86664
$ws_CRP.Range("L107").Value = 2470
$ws_CRP.Range("M107").Value = 1280.13336
$ws_CRP.Range("N107").Value = -6310
$ws_CRP.Range("H113").Value = 3580.4062
$ws_CRP.Range("I113").Value = 3048.96
$ws_CRP.Range("J113").Value = 5478.4287
$ws_CRP.Range("K113").Value = 3048.96
$ws_CRP.Range("L113").Value = 5478.4287
$ws_CRP.Range("M113").Value = -878.96
$ws_CRP.Range("N113").Value = -9818.4287
$ws_CRP.Range("H141").Value = 59731.125
$ws_CRP.Range("I141").Value = 0
$ws_CRP.Range("J141").Value = 59731.125
$ws_CRP.Range("K141").Value = 0
$ws_CRP.Range("L141").Value = 59731.125
$ws_CRP.Range("M141").ClearContents()
$ws_CRP.Range("N141").Value = -70091.125

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H58").Value = 2500
$ws_CUL.Range("I58").Value = 2500
$ws_CUL.Range("K58").Value = 7500
$ws_CUL.Range("M58").Value = -7372
$ws_CUL.Range("H112").Value = 5442.3335
$ws_CUL.Range("I112").Value = 5263.5
$ws_CUL.Range("K112").Value = 15790.5
$ws_CUL.Range("M112").Value = -14682.5
$ws_CUL.Range("H132").Value = 4295.724
$ws_CUL.Range("J132").Value = 5719.857
$ws_CUL.Range("L132").Value = 51478.713
$ws_CUL.Range("N132").Value = -56538.713
$ws_CUL.Range("H137").Value = 90357.46000000001
$ws_CUL.Range("J137").Value = 82418.21000000001
$ws_CUL.Range("L137").Value = 247254.63
$ws_CUL.Range("N137").Value = -257454.63
$ws_CUL.Range("H139").Value = 154388.7
$ws_CUL.Range("I139").Value = 233290.47
$ws_CUL.Range("K139").Value = 699871.41
$ws_CUL.Range("M139").Value = -694731.41

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H97").Value = 1928.7222
$ws_GSM.Range("I97").Value = 2093.0667
$ws_GSM.Range("K97").Value = 2093.0667
$ws_GSM.Range("M97").Value = -1597.0667
$ws_GSM.Range("H102").Value = 4588.6665
$ws_GSM.Range("I102").Value = 2112
$ws_GSM.Range("K102").Value = 2112
$ws_GSM.Range("M102").Value = -490
$ws_GSM.Range("H107").Value = 501089.56
$ws_GSM.Range("J107").Value = 2136.75
$ws_GSM.Range("L107").Value = 2136.75
$ws_GSM.Range("N107").Value = -5976.75
$ws_GSM.Range("H132").Value = 3147.1875
$ws_GSM.Range("I132").Value = 2261.3667
$ws_GSM.Range("J132").Value = 4623.5557
$ws_GSM.Range("K132").Value = 6784.1001
$ws_GSM.Range("L132").Value = 13870.6671
$ws_GSM.Range("M132").Value = -4254.1001
$ws_GSM.Range("N132").Value = -18930.6671

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 1811.125
$ws_LTW.Range("I22").Value = 1047.5
$ws_LTW.Range("J22").Value = 2574.75
$ws_LTW.Range("K22").Value = 1047.5
$ws_LTW.Range("L22").Value = 2574.75
$ws_LTW.Range("M22").Value = -752.5
$ws_LTW.Range("N22").Value = -3164.75
$ws_LTW.Range("H27").Value = 1811.125
$ws_LTW.Range("I27").Value = 1047.5
$ws_LTW.Range("J27").Value = 2574.75
$ws_LTW.Range("K27").Value = 1047.5
$ws_LTW.Range("L27").Value = 2574.75
$ws_LTW.Range("M27").Value = -940.5
$ws_LTW.Range("N27").Value = -2788.75
$ws_LTW.Range("H29").Value = 19500
$ws_LTW.Range("I29").Value = 19500
$ws_LTW.Range("K29").Value = 19500
$ws_LTW.Range("M29").Value = -19205
$ws_LTW.Range("H93").Value = 7419.778
$ws_LTW.Range("I93").Value = 8662.875
$ws_LTW.Range("J93").Value = 6425.3
$ws_LTW.Range("K93").Value = 8662.875
$ws_LTW.Range("L93").Value = 6425.3
$ws_LTW.Range("M93").Value = -7414.875
$ws_LTW.Range("N93").Value = -8921.299999999999
$ws_LTW.Range("H100").Value = 4818.273
$ws_LTW.Range("J100").Value = 6167.1665
$ws_LTW.Range("L100").Value = 6167.1665
$ws_LTW.Range("N100").Value = -7249.1665
$ws_LTW.Range("H132").Value = 12826521
$ws_LTW.Range("I132").Value = 17860380
$ws_LTW.Range("K132").Value = 53581140
$ws_LTW.Range("M132").Value = -53578610

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H100").Value = 337.8
$ws_WVR.Range("I100").Value = 337.8
$ws_WVR.Range("K100").Value = 675.6
$ws_WVR.Range("M100").Value = -134.6
$ws_WVR.Range("H107").Value = 1221.2858
$ws_WVR.Range("I107").Value = 1200.25
$ws_WVR.Range("K107").Value = 3600.75
$ws_WVR.Range("M107").Value = -1680.75
$ws_WVR.Range("H122").Value = 72657.46000000001
$ws_WVR.Range("I122").Value = 100565.05
$ws_WVR.Range("K122").Value = 301695.15
$ws_WVR.Range("M122").Value = -299245.15
$ws_WVR.Range("H132").Value = 13176086
$ws_WVR.Range("I132").Value = 21754686
$ws_WVR.Range("J132").Value = 22232.268
$ws_WVR.Range("K132").Value = 65264058
$ws_WVR.Range("L132").Value = 66696.804
$ws_WVR.Range("M132").Value = -65261528
$ws_WVR.Range("N132").Value = -71756.804

